$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 28903.367
$ws.Range("J93").Value = 28903.367
$ws.Range("L93").Value = 28903.367
$ws.Range("N93").Value = -33895.367

$ws.Range("H100").Value = 28572544
$ws.Range("I100").Value = 28572544
$ws.Range("K100").Value = 28572544
$ws.Range("M100").Value = -28572003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5168.654
$ws.Range("I32").Value = 4005.8696
$ws.Range("J32").Value = 14083.333
$ws.Range("K32").Value = 4005.8696
$ws.Range("L32").Value = 14083.333
$ws.Range("M32").Value = -3718.8696
$ws.Range("N32").Value = -14657.333

$ws.Range("H61").Value = 1667.4445
$ws.Range("I61").Value = 1320
$ws.Range("J61").Value = 2101.75
$ws.Range("K61").Value = 1320
$ws.Range("L61").Value = 2101.75
$ws.Range("M61").Value = -1108
$ws.Range("N61").Value = -2525.75

$ws.Range("H63").Value = 13854752
$ws.Range("I63").Value = 23087752
$ws.Range("J63").Value = 5250
$ws.Range("K63").Value = 23087752
$ws.Range("L63").Value = 5250
$ws.Range("M63").Value = -23087066
$ws.Range("N63").Value = -6622

$ws.Range("H66").Value = 13854752
$ws.Range("I66").Value = 23087752
$ws.Range("J66").Value = 5250
$ws.Range("K66").Value = 115438760
$ws.Range("L66").Value = 26250
$ws.Range("M66").Value = -115435328
$ws.Range("N66").Value = -33114

$ws.Range("H110").Value = 1147.7142
$ws.Range("I110").Value = 1088.2
$ws.Range("J110").Value = 1296.5
$ws.Range("K110").Value = 1088.2
$ws.Range("L110").Value = 1296.5
$ws.Range("M110").Value = 956.8
$ws.Range("N110").Value = -5386.5

$ws.Range("H136").Value = 1667.4445
$ws.Range("I136").Value = 1320
$ws.Range("J136").Value = 2101.75
$ws.Range("K136").Value = 3960
$ws.Range("L136").Value = 6305.25
$ws.Range("M136").Value = -1410
$ws.Range("N136").Value = -11405.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2414.2856
$ws.Range("I86").Value = 1633.3334
$ws.Range("K86").Value = 1633.3334
$ws.Range("M86").Value = -510.3334

$ws.Range("H89").Value = 2414.2856
$ws.Range("I89").Value = 1633.3334
$ws.Range("K89").Value = 8166.666999999999
$ws.Range("M89").Value = -2550.666999999999

$ws.Range("H94").Value = 1283.2222
$ws.Range("I94").Value = 799.8570999999999
$ws.Range("J94").Value = 2975
$ws.Range("K94").Value = 799.8570999999999
$ws.Range("L94").Value = 2975
$ws.Range("M94").Value = -348.8570999999999
$ws.Range("N94").Value = -3877

$ws.Range("H95").Value = 32475
$ws.Range("J95").Value = 32475
$ws.Range("L95").Value = 32475
$ws.Range("N95").Value = -37967

$ws.Range("H105").Value = 1740.9539
$ws.Range("I105").Value = 1745.4921
$ws.Range("K105").Value = 1745.4921
$ws.Range("M105").Value = 1.507900000000063

$ws.Range("H107").Value = 2114.5881
$ws.Range("I107").Value = 1459.3636
$ws.Range("K107").Value = 1459.3636
$ws.Range("M107").Value = 460.6364000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2460.7112
$ws.Range("I31").Value = 1094.421
$ws.Range("J31").Value = 3459.1538
$ws.Range("K31").Value = 1094.421
$ws.Range("L31").Value = 3459.1538
$ws.Range("M31").Value = -799.421
$ws.Range("N31").Value = -4049.1538

$ws.Range("H34").Value = 2460.7112
$ws.Range("I34").Value = 1094.421
$ws.Range("J34").Value = 3459.1538
$ws.Range("K34").Value = 1094.421
$ws.Range("L34").Value = 3459.1538
$ws.Range("M34").Value = -892.421
$ws.Range("N34").Value = -3863.1538

$ws.Range("H97").Value = 34999.855
$ws.Range("J97").Value = 34999.855
$ws.Range("L97").Value = 34999.855
$ws.Range("N97").Value = -36981.855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 118.875
$ws.Range("I33").Value = 110.5
$ws.Range("J33").Value = 144
$ws.Range("K33").Value = 663
$ws.Range("L33").Value = 864
$ws.Range("M33").Value = -380
$ws.Range("N33").Value = -1430

$ws.Range("H131").Value = 687.7835
$ws.Range("J131").Value = 807.2763
$ws.Range("L131").Value = 2421.8289
$ws.Range("N131").Value = -12501.8289

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 83336000
$ws.Range("I80").Value = 83336000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 83336000
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = -83335002

$ws.Range("H83").Value = 83336000
$ws.Range("I83").Value = 83336000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 416680000
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = -416675008

$ws.Range("H113").Value = 1481.3182
$ws.Range("I113").Value = 1663.3077
$ws.Range("K113").Value = 1663.3077
$ws.Range("M113").Value = 506.6922999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4293.9688
$ws.Range("I82").Value = 5615.4
$ws.Range("J82").Value = 2091.5833
$ws.Range("K82").Value = 5615.4
$ws.Range("L82").Value = 2091.5833
$ws.Range("M82").Value = -5254.4
$ws.Range("N82").Value = -2813.5833

$ws.Range("H85").Value = 4293.9688
$ws.Range("I85").Value = 5615.4
$ws.Range("J85").Value = 2091.5833
$ws.Range("K85").Value = 5615.4
$ws.Range("L85").Value = 2091.5833
$ws.Range("M85").Value = -4367.4
$ws.Range("N85").Value = -4587.5833

$ws.Range("H87").Value = 39750
$ws.Range("I87").Value = 9000
$ws.Range("J87").Value = 50000
$ws.Range("K87").Value = 9000
$ws.Range("L87").Value = 50000
$ws.Range("M87").Value = -7877
$ws.Range("N87").Value = -52246

$ws.Range("H90").Value = 39750
$ws.Range("I90").Value = 9000
$ws.Range("J90").Value = 50000
$ws.Range("K90").Value = 27000
$ws.Range("L90").Value = 150000
$ws.Range("M90").Value = -21384
$ws.Range("N90").Value = -161232

$ws.Range("H136").Value = 3896
$ws.Range("I136").Value = 1439.8572
$ws.Range("J136").Value = 7022
$ws.Range("K136").Value = 4319.571599999999
$ws.Range("L136").Value = 21066
$ws.Range("M136").Value = -1769.571599999999
$ws.Range("N136").Value = -26166

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2000
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""

$ws.Range("H84").Value = 2000
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""

$ws.Range("H107").Value = 576.6539
$ws.Range("I107").Value = 541.25
$ws.Range("J107").Value = 633.3
$ws.Range("K107").Value = 1623.75
$ws.Range("L107").Value = 1899.9
$ws.Range("M107").Value = 296.25
$ws.Range("N107").Value = -5739.9

$ws.Range("H122").Value = 3425
$ws.Range("I122").Value = 1751.7
$ws.Range("K122").Value = 5255.1
$ws.Range("M122").Value = -2805.1

$ws.Range("H136").Value = 35595.6
$ws.Range("I136").Value = 100000
$ws.Range("J136").Value = 19494.5
$ws.Range("K136").Value = 300000
$ws.Range("L136").Value = 58483.5
$ws.Range("M136").Value = -297450
$ws.Range("N136").Value = -63583.5
